$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 138
$ws.Range("F5").Value = 384
$ws.Range("F6").Value = 775
$ws.Range("F7").Value = 219
$ws.Range("F8").Value = 1095
$ws.Range("F9").Value = 286
$ws.Range("F10").Value = 1306
$ws.Range("F12").Value = 635
$ws.Range("F13").Value = 169
$ws.Range("F14").Value = 497
$ws.Range("F18").Value = 841
$ws.Range("F19").Value = 2592
$ws.Range("F27").Value = 586
$ws.Range("F28").Value = 968
$ws.Range("F31").Value = 253
$ws.Range("F32").Value = 1042
$ws.Range("F33").Value = 72

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1053
$ws.Range("F5").Value = 1053
$ws.Range("F7").Value = 16
$ws.Range("F10").Value = 321
$ws.Range("F14").Value = 592
$ws.Range("F20").Value = 609
$ws.Range("F22").Value = 37
$ws.Range("F23").Value = 3
$ws.Range("F24").Value = 295
$ws.Range("F25").Value = 265
$ws.Range("F26").Value = 3799
$ws.Range("F31").Value = 29
$ws.Range("F33").Value = 138

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1770
$ws.Range("F5").Value = 2413
$ws.Range("F6").Value = 997
$ws.Range("F9").Value = 1259
$ws.Range("F10").Value = 333

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1770
$ws.Range("F4").Value = 2413
$ws.Range("F6").Value = 997
$ws.Range("F7").Value = 1259
$ws.Range("F8").Value = 333
$ws.Range("F10").Value = 138
$ws.Range("F11").Value = 384
$ws.Range("F12").Value = 775
$ws.Range("F13").Value = 219
$ws.Range("F15").Value = 1095
$ws.Range("F16").Value = 286
$ws.Range("F18").Value = 635
$ws.Range("F19").Value = 1053
$ws.Range("F20").Value = 497
$ws.Range("F22").Value = 16
$ws.Range("F24").Value = 841
$ws.Range("F25").Value = 2592
$ws.Range("F29").Value = 321
$ws.Range("F32").Value = 586
$ws.Range("F33").Value = 968
$ws.Range("F34").Value = 592
$ws.Range("F35").Value = 592
$ws.Range("F39").Value = 253
$ws.Range("F43").Value = 37
$ws.Range("F44").Value = 295
$ws.Range("F45").Value = 295
$ws.Range("F46").Value = 265
$ws.Range("F47").Value = 1042
$ws.Range("F49").Value = 138

